# SIR framework cleanup: remove the "initprev" (Initial prevalence) parameter row
# from the Parameters sheet. This removes the now-unused shared strings
# ("initprev" / "Initial prevalence") and shifts every subsequent parameter row
# up by one, which matches FOI now depending directly/nonlinearly on prevalence
# instead of via a separate initial-prevalence parameter.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Select row 2 (the "initprev" row) before deleting it, so the sheet's
# remaining selection matches what Excel leaves behind after a row delete
# (the entire new row 2, i.e. what used to be row 3).
$ws.Rows.Item(2).Select() | Out-Null
$ws.Rows.Item(2).Delete() | Out-Null
